# Fixed Bento 80 Test scripts
# Appends "order by ... LIMIT 100" clauses to the three Cypher queries
# stored in column B (rows 2-4: Cases, Samples, Files), and adjusts the
# view's top-left cell and the affected rows' heights to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: Cases query ---
$b2 = $ws.Cells.Item(2, 2).Value()
$newB2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"
$ws.Cells.Item(2, 2).Value = $newB2

# --- Row 3: Samples query ---
$b3 = $ws.Cells.Item(3, 2).Value()
$newB3 = $b3 + "`norder By samp.sample_id ASC LIMIT 100"
$ws.Cells.Item(3, 2).Value = $newB3

# --- Row 4: Files query (replace the old trailing "order by" clause) ---
$b4 = $ws.Cells.Item(4, 2).Value()
$newB4 = $b4.Replace("`n    order by f.file_name", "`n order By f.file_name ASC LIMIT 100")
$ws.Cells.Item(4, 2).Value = $newB4

# --- Row heights grew because the wrapped text now spans one more line ---
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# --- View scrolled up by one row ---
$ws.Application.ActiveWindow.ScrollRow = 3

$wb.Save()
